$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")
Write-Host $ws.Name
Write-Host $ws.Range("F2").Value
